$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Computadores")

# Helper: write a literal, date-shaped string into a merged "A:B" cell of the
# MANTENIMIENTO PREVENTIVO table without Excel auto-converting it into a real
# date serial (which would also swap in a new number-format style).
#
#   1. Stage the literal text as a formula result ("=""2025-05-14""") - a
#      formula is never subject to the "looks like a date" auto-typing that
#      plain .Value/.Formula assignment of a bare string triggers.
#   2. Copy / PasteSpecial(values) to collapse the formula down into a plain
#      text cell. (This is also what strips the formula tag so the cell is a
#      real literal, matching the target file - but it has the side effect of
#      unmerging the A:B pair.)
#   3. Re-merge A:B.
#   4. Copy formatting from the same columns on an untouched donor row (one
#      that still carries the table's default, unedited style) and
#      PasteSpecial(formats) it back onto the row, since step 3 resets the
#      merged range's style.
function Set-LiteralDateText {
    param($ws, [string]$cellRef, [string]$text, [string]$donorRowRange)

    # Capture the merged extent (e.g. "A34:B34") *before* anything below
    # collapses the merge, so it can be restored afterwards.
    $mergeRef = $ws.Range($cellRef).MergeArea.Address($false, $false)

    $ws.Range($cellRef).Formula = '="' + $text + '"'
    $ws.Range($cellRef).Copy() | Out-Null
    $ws.Range($cellRef).PasteSpecial(-4163) | Out-Null   # xlPasteValues
    $excel.CutCopyMode = 0

    $ws.Range($mergeRef).Merge() | Out-Null

    $ws.Range($donorRowRange).Copy() | Out-Null
    $ws.Range($mergeRef).PasteSpecial(-4122) | Out-Null  # xlPasteFormats
    $excel.CutCopyMode = 0
}

# --- Row 34: fecha, actividad, funcionario TIC ---
Set-LiteralDateText $ws "A34" "2025-05-14" "A36:B36"
$ws.Range("C34").Value = "uuuuuuu"
$ws.Range("M34").Value = "SV Romero Romero Miguel Ángel"

# --- Row 35: fecha, actividad, funcionario TIC ---
Set-LiteralDateText $ws "A35" "2025-05-26" "A36:B36"
$ws.Range("C35").Value = "clon 2222222"
$ws.Range("M35").Value = "SV Romero Romero Miguel Ángel"
